# Tomorrow gas prices - Canadian cities: refresh regular/premium/diesel
# prices and day-over-day change indicators from the latest CSV pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $text) {
    # Leading apostrophe forces Excel to keep numeric-looking
    # strings (e.g. "145.9") stored as text, matching columns
    # B/C/D/J/L which hold inline-string values, not numbers.
    $sheet.Range($addr).Value = "'" + $text
}

# Row 2
Set-TextCell $ws "B2" '145.9'
Set-TextCell $ws "C2" '175.9 ⮟ -1'
Set-TextCell $ws "D2" '210.9 n/c'
Set-TextCell $ws "J2" '⮟'
Set-TextCell $ws "L2" '⮟ -1'
$ws.Range("I2").Value = 1.46
$ws.Range("K2").Value = -1

# Row 3
Set-TextCell $ws "B3" '161.9'
Set-TextCell $ws "C3" '191.9 ⮟ -1'
Set-TextCell $ws "D3" '230.9 n/c'
Set-TextCell $ws "J3" '⮟'
Set-TextCell $ws "L3" '⮟ -1'
$ws.Range("I3").Value = 1.62
$ws.Range("K3").Value = -1

# Row 4
Set-TextCell $ws "B4" '170.9'
Set-TextCell $ws "C4" '199.9 n/c'
Set-TextCell $ws "D4" '229.9 n/c'
$ws.Range("I4").Value = 1.71

# Row 5
Set-TextCell $ws "B5" '118.9'
Set-TextCell $ws "C5" '142.9 ⮝ 7'
Set-TextCell $ws "D5" '176.9 n/c'
Set-TextCell $ws "J5" '⮝'
Set-TextCell $ws "L5" '⮝ 7'
$ws.Range("I5").Value = 1.19
$ws.Range("K5").Value = 7

# Row 6
Set-TextCell $ws "B6" '144.9'
Set-TextCell $ws "C6" '174.9 ⮟ -1'
Set-TextCell $ws "D6" '210.9 n/c'
Set-TextCell $ws "J6" '⮟'
Set-TextCell $ws "L6" '⮟ -1'
$ws.Range("I6").Value = 1.45
$ws.Range("K6").Value = -1

# Row 7
Set-TextCell $ws "B7" '145.9'
Set-TextCell $ws "C7" '175.9 ⮟ -1'
Set-TextCell $ws "D7" '210.9 n/c'
Set-TextCell $ws "J7" '⮟'
Set-TextCell $ws "L7" '⮟ -1'
$ws.Range("I7").Value = 1.46
$ws.Range("K7").Value = -1

# Row 8
Set-TextCell $ws "B8" '156.5'
Set-TextCell $ws "C8" '165.6 ⮝ 6'
Set-TextCell $ws "D8" '226.4 ⮝ 3'
Set-TextCell $ws "J8" '⮝'
Set-TextCell $ws "L8" '⮝ 6'
$ws.Range("I8").Value = 1.56
$ws.Range("K8").Value = 6

# Row 9
Set-TextCell $ws "B9" '145.9'
Set-TextCell $ws "C9" '165.9 ⮟ -1'
Set-TextCell $ws "D9" '214.9 n/c'
Set-TextCell $ws "J9" '⮟'
Set-TextCell $ws "L9" '⮟ -1'
$ws.Range("I9").Value = 1.46
$ws.Range("K9").Value = -1

# Row 10
Set-TextCell $ws "B10" '116.9'
Set-TextCell $ws "C10" '133.9 n/c'
Set-TextCell $ws "D10" '173.9 n/c'
$ws.Range("I10").Value = 1.17

# Row 11
Set-TextCell $ws "B11" '154.4'
Set-TextCell $ws "C11" '161.7 n/c'
Set-TextCell $ws "D11" '209.2 n/c'
$ws.Range("I11").Value = 1.54

# Row 12
Set-TextCell $ws "B12" '145.9'
Set-TextCell $ws "C12" '175.9 ⮟ -1'
Set-TextCell $ws "D12" '210.9 n/c'
Set-TextCell $ws "J12" '⮟'
Set-TextCell $ws "L12" '⮟ -1'
$ws.Range("I12").Value = 1.46
$ws.Range("K12").Value = -1

# Row 13
Set-TextCell $ws "B13" '147.7'
Set-TextCell $ws "C13" '157.9 ⮝ 7'
Set-TextCell $ws "D13" '203.0 ⮝ 3'
Set-TextCell $ws "J13" '⮝'
Set-TextCell $ws "L13" '⮝ 6'
$ws.Range("I13").Value = 1.48
$ws.Range("K13").Value = 6

# Row 14
Set-TextCell $ws "B14" '145.9'
Set-TextCell $ws "C14" '175.9 ⮟ -1'
Set-TextCell $ws "D14" '210.9 n/c'
Set-TextCell $ws "J14" '⮟'
Set-TextCell $ws "L14" '⮟ -1'
$ws.Range("I14").Value = 1.46
$ws.Range("K14").Value = -1

# Row 15
Set-TextCell $ws "B15" '152.9'
Set-TextCell $ws "C15" '181.9 n/c'
Set-TextCell $ws "D15" '200.9 n/c'
$ws.Range("I15").Value = 1.53

# Row 16
Set-TextCell $ws "B16" '145.9'
Set-TextCell $ws "C16" '167.9 n/c'
Set-TextCell $ws "D16" '200.9 n/c'
$ws.Range("I16").Value = 1.46

# Row 17
Set-TextCell $ws "B17" '145.9'
Set-TextCell $ws "C17" '167.9 ⮟ -1'
Set-TextCell $ws "D17" '216.9 n/c'
Set-TextCell $ws "J17" '⮟'
Set-TextCell $ws "L17" '⮟ -1'
$ws.Range("I17").Value = 1.46
$ws.Range("K17").Value = -1

# Row 18
Set-TextCell $ws "B18" '144.9'
Set-TextCell $ws "C18" '174.9 ⮟ -1'
Set-TextCell $ws "D18" '210.9 n/c'
Set-TextCell $ws "J18" '⮟'
Set-TextCell $ws "L18" '⮟ -1'
$ws.Range("I18").Value = 1.45
$ws.Range("K18").Value = -1

# Row 19
Set-TextCell $ws "B19" '145.9'
Set-TextCell $ws "C19" '175.9 ⮟ -1'
Set-TextCell $ws "D19" '210.9 n/c'
Set-TextCell $ws "J19" '⮟'
Set-TextCell $ws "L19" '⮟ -1'
$ws.Range("I19").Value = 1.46
$ws.Range("K19").Value = -1

# Row 20
Set-TextCell $ws "B20" '145.9'
Set-TextCell $ws "C20" '175.9 ⮟ -1'
Set-TextCell $ws "D20" '210.9 n/c'
Set-TextCell $ws "J20" '⮟'
Set-TextCell $ws "L20" '⮟ -1'
$ws.Range("I20").Value = 1.46
$ws.Range("K20").Value = -1

# Row 21
Set-TextCell $ws "B21" '152.9'
Set-TextCell $ws "C21" '161.8 n/c'
Set-TextCell $ws "D21" '209.2 n/c'
$ws.Range("I21").Value = 1.53

# Row 22
Set-TextCell $ws "B22" '145.9'
Set-TextCell $ws "C22" '175.9 ⮟ -1'
Set-TextCell $ws "D22" '210.9 n/c'
Set-TextCell $ws "J22" '⮟'
Set-TextCell $ws "L22" '⮟ -1'
$ws.Range("I22").Value = 1.46
$ws.Range("K22").Value = -1

# Row 23
Set-TextCell $ws "B23" '145.9'
Set-TextCell $ws "C23" '175.9 ⮟ -1'
Set-TextCell $ws "D23" '210.9 n/c'
Set-TextCell $ws "J23" '⮟'
Set-TextCell $ws "L23" '⮟ -1'
$ws.Range("I23").Value = 1.46
$ws.Range("K23").Value = -1

# Row 24
Set-TextCell $ws "B24" '144.9'
Set-TextCell $ws "C24" '174.9 ⮟ -1'
Set-TextCell $ws "D24" '210.9 n/c'
Set-TextCell $ws "J24" '⮟'
Set-TextCell $ws "L24" '⮟ -1'
$ws.Range("I24").Value = 1.45
$ws.Range("K24").Value = -1

# Row 25
Set-TextCell $ws "B25" '145.9'
Set-TextCell $ws "C25" '175.9 ⮟ -1'
Set-TextCell $ws "D25" '222.9 n/c'
Set-TextCell $ws "J25" '⮟'
Set-TextCell $ws "L25" '⮟ -1'
$ws.Range("I25").Value = 1.46
$ws.Range("K25").Value = -1

# Row 26
Set-TextCell $ws "C26" '154.9 ⮟ -1'
Set-TextCell $ws "D26" '211.9 n/c'
Set-TextCell $ws "J26" '⮟'
Set-TextCell $ws "L26" '⮟ -1'
$ws.Range("K26").Value = -1

# Row 27
Set-TextCell $ws "B27" '152.9'
Set-TextCell $ws "C27" '174.9 n/c'
Set-TextCell $ws "D27" '188.9 n/c'
$ws.Range("I27").Value = 1.53

# Row 28
Set-TextCell $ws "B28" '160.9'
Set-TextCell $ws "C28" '181.9 ⮟ -1'
Set-TextCell $ws "D28" '234.9 n/c'
Set-TextCell $ws "J28" '⮟'
Set-TextCell $ws "L28" '⮟ -1'
$ws.Range("I28").Value = 1.61
$ws.Range("K28").Value = -1

# Row 29
Set-TextCell $ws "B29" '137.9'
Set-TextCell $ws "C29" '161.9 n/c'
Set-TextCell $ws "D29" '189.9 n/c'
$ws.Range("I29").Value = 1.38

# Row 30
Set-TextCell $ws "B30" '136.9'
Set-TextCell $ws "C30" '157.9 n/c'
$ws.Range("I30").Value = 1.37

# Row 31
Set-TextCell $ws "B31" '145.9'
Set-TextCell $ws "C31" '175.9 ⮟ -1'
Set-TextCell $ws "D31" '210.9 n/c'
Set-TextCell $ws "J31" '⮟'
Set-TextCell $ws "L31" '⮟ -1'
$ws.Range("I31").Value = 1.46
$ws.Range("K31").Value = -1

# Row 32
Set-TextCell $ws "B32" '154.6'
Set-TextCell $ws "C32" '162.4 n/c'
Set-TextCell $ws "D32" '209.2 n/c'
$ws.Range("I32").Value = 1.55

# Row 33
Set-TextCell $ws "B33" '162.6'
Set-TextCell $ws "D33" '223.9 n/c'

# Row 34
Set-TextCell $ws "B34" '129.9'
Set-TextCell $ws "C34" '178.9 n/c'
Set-TextCell $ws "D34" '212.9 n/c'
$ws.Range("I34").Value = 1.3

# Row 35
Set-TextCell $ws "C35" '182.9 n/c'
Set-TextCell $ws "D35" '217.9 n/c'

# Row 36
Set-TextCell $ws "B36" '155.9'
Set-TextCell $ws "C36" '179.9 n/c'
Set-TextCell $ws "D36" '185.9 n/c'
$ws.Range("I36").Value = 1.56

# Row 37
Set-TextCell $ws "B37" '145.9'
Set-TextCell $ws "C37" '175.9 ⮟ -1'
Set-TextCell $ws "D37" '210.9 n/c'
Set-TextCell $ws "J37" '⮟'
Set-TextCell $ws "L37" '⮟ -1'
$ws.Range("I37").Value = 1.46
$ws.Range("K37").Value = -1

# Row 38
Set-TextCell $ws "B38" '145.9'
Set-TextCell $ws "C38" '175.9 ⮟ -1'
Set-TextCell $ws "D38" '210.9 n/c'
Set-TextCell $ws "J38" '⮟'
Set-TextCell $ws "L38" '⮟ -1'
$ws.Range("I38").Value = 1.46
$ws.Range("K38").Value = -1

# Row 39
Set-TextCell $ws "B39" '147.9'
Set-TextCell $ws "C39" '165.9 n/c'
$ws.Range("I39").Value = 1.48

